$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2021-12-31) was added to the data set.
# It is inserted as a new row 20, pushing all existing rows from 20
# downward by one (old row 20 becomes row 21, ..., old row 126 becomes
# row 127). This matches the observed change in dimension from
# A1:R126 to A1:R127.
$ws.Rows("20:20").Insert()

$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value = "12/31/2021"
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = 100112003
$ws.Cells.Item(20, 7).Value = "Ajo"
$ws.Cells.Item(20, 8).Value = "Chino"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 270
$ws.Cells.Item(20, 11).Value = 16000
$ws.Cells.Item(20, 12).Value = 17000
$ws.Cells.Item(20, 13).Value = 16444
$ws.Cells.Item(20, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(20, 15).Value = "China"
$ws.Cells.Item(20, 16).Value = 1644
$ws.Cells.Item(20, 17).Value = 10
$ws.Cells.Item(20, 18).Value = "Hortaliza"
